$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "Category" column before the existing "Expense" column (B),
# shifting the old Expense column to C and the old Amount column to D.
$ws.Columns("B").EntireColumn.Insert() | Out-Null

# New column header
$ws.Range("B1").Value = "Category"

# Populate the new Category column with a broader category per expense,
# and make the old "Expense" column (now C) carry more specific labels.
$ws.Range("B2").Value = "Food"
$ws.Range("C2").Value = "Jersey Mike's"

$ws.Range("B3").Value = "Food"
$ws.Range("C3").Value = "Food"

$ws.Range("B4").Value = "Gas"
$ws.Range("C4").Value = "Gas"

$ws.Range("B5").Value = "Bill"
$ws.Range("C5").Value = "Water Bill"

$ws.Range("B6").Value = "Bill"
$ws.Range("C6").Value = "Electric Bill"

# Match the column width Excel would apply to the newly inserted column
# (mirrors column A's width, same as the source column it was inserted next to).
$ws.Columns("B").ColumnWidth = $ws.Columns("A").ColumnWidth

# Leave the workbook selection where the author ended up.
$ws.Range("K10").Select() | Out-Null
